$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.892066666666667
$ws.Range("H2").Value = 23.6762
$ws.Range("I2").Value = 0.1739002798877711
$ws.Range("J2").Value = 0.1739002798877711
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1318966666666667
$ws.Range("N2").Value = 0.39569
$ws.Range("O2").Value = 0.007203585535592019
$ws.Range("P2").Value = 0.00720358553559202
$ws.Range("Q2").Value = 1.040937286444445
$ws.Range("R2").Value = 9.368435578
$ws.Range("S2").Value = 0.001252705540834951
$ws.Range("T2").Value = 0.001252705540834951
$ws.Range("G3").Value = 7.892066666666667
$ws.Range("H3").Value = 23.6762
$ws.Range("I3").Value = 0.1739002798877711
$ws.Range("J3").Value = 0.1739002798877711
$ws.Range("O3").Value = 0.05718145720730849
$ws.Range("P3").Value = 0.05718145720730849
$ws.Range("Q3").Value = 8.262872788311112
$ws.Range("R3").Value = 74.3658550948
$ws.Range("S3").Value = 0.009943871412741551
$ws.Range("T3").Value = 0.009943871412741551
$ws.Range("G4").Value = 7.892066666666667
$ws.Range("H4").Value = 23.6762
$ws.Range("I4").Value = 0.1739002798877711
$ws.Range("J4").Value = 0.1739002798877711
$ws.Range("M4").Value = 3.572802
$ws.Range("N4").Value = 10.718406
$ws.Range("O4").Value = 0.1951299108549691
$ws.Range("P4").Value = 0.1951299108549691
$ws.Range("Q4").Value = 28.1967915708
$ws.Range("R4").Value = 253.7711241372
$ws.Range("S4").Value = 0.03393314611215494
$ws.Range("T4").Value = 0.03393314611215494
$ws.Range("G5").Value = 7.892066666666667
$ws.Range("H5").Value = 23.6762
$ws.Range("I5").Value = 0.1739002798877711
$ws.Range("J5").Value = 0.1739002798877711
$ws.Range("M5").Value = 13.55818
$ws.Range("N5").Value = 40.67453999999999
$ws.Range("O5").Value = 0.7404850464021304
$ws.Range("P5").Value = 0.7404850464021304
$ws.Range("Q5").Value = 107.0020604386667
$ws.Range("R5").Value = 963.0185439479999
$ws.Range("S5").Value = 0.1287705568220396
$ws.Range("T5").Value = 0.1287705568220396
$ws.Range("I6").Value = 0.3815924715300191
$ws.Range("J6").Value = 0.3815924715300191
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.1318966666666667
$ws.Range("N6").Value = 0.39569
$ws.Range("O6").Value = 0.007203585535592019
$ws.Range("P6").Value = 0.00720358553559202
$ws.Range("Q6").Value = 2.284147168126667
$ws.Range("R6").Value = 20.55732451314
$ws.Range("S6").Value = 0.002748834008404455
$ws.Range("T6").Value = 0.002748834008404455
$ws.Range("I7").Value = 0.3815924715300191
$ws.Range("J7").Value = 0.3815924715300191
$ws.Range("O7").Value = 0.05718145720730849
$ws.Range("P7").Value = 0.05718145720730849
$ws.Range("S7").Value = 0.02182001358142487
$ws.Range("T7").Value = 0.02182001358142487
$ws.Range("I8").Value = 0.3815924715300191
$ws.Range("J8").Value = 0.3815924715300191
$ws.Range("M8").Value = 3.572802
$ws.Range("N8").Value = 10.718406
$ws.Range("O8").Value = 0.1951299108549691
$ws.Range("P8").Value = 0.1951299108549691
$ws.Range("Q8").Value = 61.872720341004
$ws.Range("R8").Value = 556.8544830690361
$ws.Range("S8").Value = 0.07446010495257994
$ws.Range("T8").Value = 0.07446010495257996
$ws.Range("I9").Value = 0.3815924715300191
$ws.Range("J9").Value = 0.3815924715300191
$ws.Range("M9").Value = 13.55818
$ws.Range("N9").Value = 40.67453999999999
$ws.Range("O9").Value = 0.7404850464021304
$ws.Range("P9").Value = 0.7404850464021304
$ws.Range("Q9").Value = 234.79652090236
$ws.Range("R9").Value = 2113.16868812124
$ws.Range("S9").Value = 0.2825635189876098
$ws.Range("T9").Value = 0.2825635189876098
$ws.Range("G10").Value = 7.716272666666666
$ws.Range("H10").Value = 23.148818
$ws.Range("I10").Value = 0.1700266904854272
$ws.Range("J10").Value = 0.1700266904854272
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.1318966666666667
$ws.Range("N10").Value = 0.39569
$ws.Range("O10").Value = 0.007203585535592019
$ws.Range("P10").Value = 0.00720358553559202
$ws.Range("Q10").Value = 1.017750643824444
$ws.Range("R10").Value = 9.159755794419999
$ws.Range("S10").Value = 0.001224801808245405
$ws.Range("T10").Value = 0.001224801808245405
$ws.Range("G11").Value = 7.716272666666666
$ws.Range("H11").Value = 23.148818
$ws.Range("I11").Value = 0.1700266904854272
$ws.Range("J11").Value = 0.1700266904854272
$ws.Range("O11").Value = 0.05718145720730849
$ws.Range("P11").Value = 0.05718145720730849
$ws.Range("Q11").Value = 8.078819165819111
$ws.Range("R11").Value = 72.709372492372
$ws.Range("S11").Value = 0.009722373926092742
$ws.Range("T11").Value = 0.009722373926092742
$ws.Range("G12").Value = 7.716272666666666
$ws.Range("H12").Value = 23.148818
$ws.Range("I12").Value = 0.1700266904854272
$ws.Range("J12").Value = 0.1700266904854272
$ws.Range("M12").Value = 3.572802
$ws.Range("N12").Value = 10.718406
$ws.Range("O12").Value = 0.1951299108549691
$ws.Range("P12").Value = 0.1951299108549691
$ws.Range("Q12").Value = 27.568714416012
$ws.Range("R12").Value = 248.118429744108
$ws.Range("S12").Value = 0.03317729295738683
$ws.Range("T12").Value = 0.03317729295738683
$ws.Range("G13").Value = 7.716272666666666
$ws.Range("H13").Value = 23.148818
$ws.Range("I13").Value = 0.1700266904854272
$ws.Range("J13").Value = 0.1700266904854272
$ws.Range("M13").Value = 13.55818
$ws.Range("N13").Value = 40.67453999999999
$ws.Range("O13").Value = 0.7404850464021304
$ws.Range("P13").Value = 0.7404850464021304
$ws.Range("Q13").Value = 104.6186137437466
$ws.Range("R13").Value = 941.5675236937197
$ws.Range("S13").Value = 0.1259022217937022
$ws.Range("T13").Value = 0.1259022217937023
$ws.Range("G14").Value = 12.45667266666667
$ws.Range("H14").Value = 37.370018
$ws.Range("I14").Value = 0.2744805580967825
$ws.Range("J14").Value = 0.2744805580967826
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.1318966666666667
$ws.Range("N14").Value = 0.39569
$ws.Range("O14").Value = 0.007203585535592019
$ws.Range("P14").Value = 0.00720358553559202
$ws.Range("Q14").Value = 1.642993602491111
$ws.Range("R14").Value = 14.78694242242
$ws.Range("S14").Value = 0.001977244178107207
$ws.Range("T14").Value = 0.001977244178107208
$ws.Range("G15").Value = 12.45667266666667
$ws.Range("H15").Value = 37.370018
$ws.Range("I15").Value = 0.2744805580967825
$ws.Range("J15").Value = 0.2744805580967826
$ws.Range("O15").Value = 0.05718145720730849
$ws.Range("P15").Value = 0.05718145720730849
$ws.Range("Q15").Value = 13.04194527968578
$ws.Range("R15").Value = 117.377507517172
$ws.Range("S15").Value = 0.01569519828704932
$ws.Range("T15").Value = 0.01569519828704932
$ws.Range("G16").Value = 12.45667266666667
$ws.Range("H16").Value = 37.370018
$ws.Range("I16").Value = 0.2744805580967825
$ws.Range("J16").Value = 0.2744805580967826
$ws.Range("M16").Value = 3.572802
$ws.Range("N16").Value = 10.718406
$ws.Range("O16").Value = 0.1951299108549691
$ws.Range("P16").Value = 0.1951299108549691
$ws.Range("Q16").Value = 44.50522501681201
$ws.Range("R16").Value = 400.547025151308
$ws.Range("S16").Value = 0.05355936683284733
$ws.Range("T16").Value = 0.05355936683284734
$ws.Range("G17").Value = 12.45667266666667
$ws.Range("H17").Value = 37.370018
$ws.Range("I17").Value = 0.2744805580967825
$ws.Range("J17").Value = 0.2744805580967826
$ws.Range("M17").Value = 13.55818
$ws.Range("N17").Value = 40.67453999999999
$ws.Range("O17").Value = 0.7404850464021304
$ws.Range("P17").Value = 0.7404850464021304
$ws.Range("Q17").Value = 168.8898102157466
$ws.Range("R17").Value = 1520.00829194172
$ws.Range("S17").Value = 0.2032487487987787
$ws.Range("T17").Value = 0.2032487487987787
